$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 161, shifting existing rows 161:173 down to 162:174
$ws.Rows.Item(161).Insert()

# Fill in the new row 161 with the new weekly data point
$ws.Cells.Item(161, 1).Value = 7
$ws.Cells.Item(161, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(161, 3).Value = "Ñuble"
$ws.Cells.Item(161, 4).Value = 44461
$ws.Cells.Item(161, 4).NumberFormat = $ws.Cells.Item(162, 4).NumberFormat
$ws.Cells.Item(161, 5).Value = 16
$ws.Cells.Item(161, 6).Value = 100114013
$ws.Cells.Item(161, 7).Value = "Zanahoria"
$ws.Cells.Item(161, 8).Value = "Sin especificar"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 300
$ws.Cells.Item(161, 11).Value = 8500
$ws.Cells.Item(161, 12).Value = 9000
$ws.Cells.Item(161, 13).Value = 8750
$ws.Cells.Item(161, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(161, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(161, 16).Value = 438
$ws.Cells.Item(161, 17).Value = 20
$ws.Cells.Item(161, 18).Value = "Hortaliza"
